$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1175.8334
$ws.Range("I28").Value = 346.66666
$ws.Range("J28").Value = 3663.3333
$ws.Range("K28").Value = 346.66666
$ws.Range("L28").Value = 3663.3333
$ws.Range("M28").Value = 138.33334
$ws.Range("N28").Value = -4633.3333
$ws.Range("H43").Value = 3880.2
$ws.Range("I43").Value = 6000.5
$ws.Range("J43").Value = 2466.6667
$ws.Range("K43").Value = 6000.5
$ws.Range("L43").Value = 2466.6667
$ws.Range("M43").Value = -5931.5
$ws.Range("N43").Value = -2604.6667
$ws.Range("H47").Value = 13999.667
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").Value = ""
$ws.Range("H53").Value = 860.2222
$ws.Range("I53").Value = 656.6667
$ws.Range("J53").Value = 1267.3334
$ws.Range("K53").Value = 656.6667
$ws.Range("L53").Value = 1267.3334
$ws.Range("M53").Value = -19.66669999999999
$ws.Range("N53").Value = -2541.3334
$ws.Range("H55").Value = 153.82353
$ws.Range("I55").Value = 38.11111
$ws.Range("J55").Value = 284
$ws.Range("K55").Value = 38.11111
$ws.Range("L55").Value = 284
$ws.Range("M55").Value = 175.88889
$ws.Range("N55").Value = -712
$ws.Range("H61").Value = 25148.428
$ws.Range("I61").Value = 576.25
$ws.Range("J61").Value = 57911.332
$ws.Range("K61").Value = 1728.75
$ws.Range("L61").Value = 173733.996
$ws.Range("M61").Value = -1556.75
$ws.Range("N61").Value = -174077.996
$ws.Range("H132").Value = 4172067.2
$ws.Range("I132").Value = 4635000.5
$ws.Range("K132").Value = 13905001.5
$ws.Range("M132").Value = -13902471.5
$ws.Range("H137").Value = 1772.28
$ws.Range("I137").Value = 1406.6875
$ws.Range("K137").Value = 4220.0625
$ws.Range("M137").Value = -1670.0625
$ws.Range("H141").Value = 2097.9697
$ws.Range("I141").Value = 1695.9584
$ws.Range("J141").Value = 3170
$ws.Range("K141").Value = 5087.8752
$ws.Range("L141").Value = 9510
$ws.Range("M141").Value = 92.1247999999996
$ws.Range("N141").Value = -19870
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 137
$ws.Range("I4").Value = 137
$ws.Range("K4").Value = 137
$ws.Range("M4").Value = -21
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").Value = ""
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").Value = ""
$ws.Range("H32").Value = 3875.99
$ws.Range("I32").Value = 3398.705
$ws.Range("K32").Value = 3398.705
$ws.Range("M32").Value = -3111.705
$ws.Range("H37").Value = 15000
$ws.Range("J37").Value = 15000
$ws.Range("L37").Value = 15000
$ws.Range("N37").Value = -15546
$ws.Range("H55").Value = 14764.444
$ws.Range("J55").Value = 15485
$ws.Range("L55").Value = 15485
$ws.Range("N55").Value = -16115
$ws.Range("H61").Value = 1728.4688
$ws.Range("I61").Value = 1353.8889
$ws.Range("J61").Value = 2210.0715
$ws.Range("K61").Value = 1353.8889
$ws.Range("L61").Value = 2210.0715
$ws.Range("M61").Value = -1141.8889
$ws.Range("N61").Value = -2634.0715
$ws.Range("H80").Value = 24497.5
$ws.Range("J80").Value = 24906.363
$ws.Range("L80").Value = 24906.363
$ws.Range("N80").Value = -26902.363
$ws.Range("H83").Value = 24497.5
$ws.Range("J83").Value = 24906.363
$ws.Range("L83").Value = 74719.08900000001
$ws.Range("N83").Value = -84703.08900000001
$ws.Range("H136").Value = 1728.4688
$ws.Range("I136").Value = 1353.8889
$ws.Range("J136").Value = 2210.0715
$ws.Range("K136").Value = 4061.6667
$ws.Range("L136").Value = 6630.2145
$ws.Range("M136").Value = -1511.6667
$ws.Range("N136").Value = -11730.2145
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 225
$ws.Range("I22").Value = 225
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 225
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -52
$ws.Range("N22").Value = ""
$ws.Range("H134").Value = 2571.415
$ws.Range("I134").Value = 2445.7334
$ws.Range("J134").Value = 3278.375
$ws.Range("K134").Value = 7337.2002
$ws.Range("L134").Value = 9835.125
$ws.Range("M134").Value = -4802.2002
$ws.Range("N134").Value = -14905.125
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 372.875
$ws.Range("I22").Value = 220.75
$ws.Range("J22").Value = 525
$ws.Range("K22").Value = 220.75
$ws.Range("L22").Value = 525
$ws.Range("M22").Value = 129.25
$ws.Range("N22").Value = -1225
$ws.Range("H31").Value = 42960.973
$ws.Range("I31").Value = 1019.7647
$ws.Range("K31").Value = 1019.7647
$ws.Range("M31").Value = -724.7646999999999
$ws.Range("H34").Value = 42960.973
$ws.Range("I34").Value = 1019.7647
$ws.Range("K34").Value = 1019.7647
$ws.Range("M34").Value = -817.7646999999999
$ws.Range("H58").Value = 1951.48
$ws.Range("I58").Value = 1999.5333
$ws.Range("J58").Value = 1879.4
$ws.Range("K58").Value = 1999.5333
$ws.Range("L58").Value = 1879.4
$ws.Range("M58").Value = -1796.5333
$ws.Range("N58").Value = -2285.4
$ws.Range("H134").Value = 971.8461
$ws.Range("I134").Value = 524.6061
$ws.Range("J134").Value = 3431.6667
$ws.Range("K134").Value = 1573.8183
$ws.Range("L134").Value = 10295.0001
$ws.Range("M134").Value = 961.1817000000001
$ws.Range("N134").Value = -15365.0001
$ws.Range("H136").Value = 1951.48
$ws.Range("I136").Value = 1999.5333
$ws.Range("J136").Value = 1879.4
$ws.Range("K136").Value = 5998.5999
$ws.Range("L136").Value = 5638.200000000001
$ws.Range("M136").Value = -3448.5999
$ws.Range("N136").Value = -10738.2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H127").Value = 1088.8
$ws.Range("J127").Value = 1117
$ws.Range("L127").Value = 3351
$ws.Range("N127").Value = -13271
$ws.Range("H129").Value = 161781.78
$ws.Range("J129").Value = 204670.48
$ws.Range("L129").Value = 614011.4400000001
$ws.Range("N129").Value = -624011.4400000001
$ws.Range("H131").Value = 539347.0600000001
$ws.Range("I131").Value = 763.35297
$ws.Range("J131").Value = 658255.1
$ws.Range("K131").Value = 2290.05891
$ws.Range("L131").Value = 1974765.3
$ws.Range("M131").Value = 2749.94109
$ws.Range("N131").Value = -1984845.3
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 80531.81
$ws.Range("I70").Value = 110912.69
$ws.Range("K70").Value = 110912.69
$ws.Range("M70").Value = -110642.69
$ws.Range("H73").Value = 80531.81
$ws.Range("I73").Value = 110912.69
$ws.Range("K73").Value = 110912.69
$ws.Range("M73").Value = -109976.69
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 6279.857
$ws.Range("I38").Value = 2900
$ws.Range("J38").Value = 6843.1665
$ws.Range("K38").Value = 2900
$ws.Range("L38").Value = 6843.1665
$ws.Range("M38").Value = -2427
$ws.Range("N38").Value = -7789.1665
$ws.Range("H46").Value = 43997.25
$ws.Range("J46").Value = 43997.25
$ws.Range("L46").Value = 43997.25
$ws.Range("N46").Value = -44459.25
$ws.Range("H134").Value = 43997.25
$ws.Range("J134").Value = 43997.25
$ws.Range("L134").Value = 131991.75
$ws.Range("N134").Value = -137061.75
$ws.Range("H136").Value = 14997.099
$ws.Range("I136").Value = 24428.453
$ws.Range("J136").Value = 4840.2563
$ws.Range("K136").Value = 73285.359
$ws.Range("L136").Value = 14520.7689
$ws.Range("M136").Value = -70735.359
$ws.Range("N136").Value = -19620.7689

Write-Output "done"